# Upgrade Workday Closing Email to Professional Template
# Replaces the 6-ticket planning sample data with a 4-ticket "Soporte"
# batch (ticket_id 100-103) at a single address, drops the two trailing
# rows, clears the now-unused "Accesorios" cells for the normal-priority
# rows, and switches the date column's number format from m/d/yyyy-style
# (numFmtId 14) to d-mmm (numFmtId 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the last two sample rows (former tickets "23" and "24") -- the
# new data set only has 4 rows.
$ws.Rows("6:7").Delete()

# Row 2 - urgente Soporte ticket for Juan Parez
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = "urgente"
$ws.Range("D2").Value = "Soporte"
$ws.Range("E2").Value = "Antena GPS, Corta Corriente Remoto"
$ws.Range("F2").Value = "AMERICO VESPUCIO NORTE 2341, PUDAHUEL"
$ws.Range("G2").Value = "PUDAHUEL"
$ws.Range("I2").Value = "Juan Parez"
$ws.Range("J2").Value = "PRSV10"
$ws.Range("K2").Value = "CIAL_ALIMENTOS"

# Row 3 - urgente Soporte ticket for Juan Parez
$ws.Range("B3").Value = 101
$ws.Range("C3").Value = "urgente"
$ws.Range("D3").Value = "Soporte"
$ws.Range("E3").Value = "Antena GPS"
$ws.Range("F3").Value = "AMERICO VESPUCIO NORTE 2341, PUDAHUEL"
$ws.Range("G3").Value = "PUDAHUEL"
$ws.Range("I3").Value = "Juan Parez"
$ws.Range("J3").Value = "VE041-LCSH40"
$ws.Range("K3").Value = "SOPROLE"

# Row 4 - normal Soporte ticket for Pedro Pascal (no accessories)
$ws.Range("B4").Value = 102
$ws.Range("C4").Value = "normal"
$ws.Range("D4").Value = "Soporte"
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = "AMERICO VESPUCIO NORTE 2341, PUDAHUEL"
$ws.Range("G4").Value = "PUDAHUEL"
$ws.Range("I4").Value = "Pedro Pascal"
$ws.Range("J4").Value = "VCVF53"
$ws.Range("K4").Value = "AGUNSA"

# Row 5 - normal Soporte ticket for Pedro Pascal (no accessories)
$ws.Range("B5").Value = 103
$ws.Range("C5").Value = "normal"
$ws.Range("D5").Value = "Soporte"
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = "AMERICO VESPUCIO NORTE 2341, PUDAHUEL"
$ws.Range("G5").Value = "PUDAHUEL"
$ws.Range("I5").Value = "Pedro Pascal"
$ws.Range("J5").Value = "HPWV83"
$ws.Range("K5").Value = "BRINKS"

# Date column now displays as "d-mmm" instead of the previous date format.
$ws.Range("A2:A5").NumberFormat = "d-mmm"

# Leave the cursor where the author ended up after entering the data.
$ws.Range("K6").Select()
